$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before the current row 46, shifting the existing
# rows 46:115 down to 49:118 (matches the new dimension A1:R118).
$ws.Rows("46:48").Insert()

# Row 46 (new)
$ws.Range("A46").Value = 1
$ws.Range("B46").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C46").Value = "Arica y Parinacota"
$ws.Range("D46").Value = 44915
$ws.Range("E46").Value = 15
$ws.Range("F46").Value = 100114001
$ws.Range("G46").Value = "Papa"
$ws.Range("H46").Value = "Asterix"
$ws.Range("I46").Value = "1a (cosecha)"
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 17000
$ws.Range("L46").Value = 19000
$ws.Range("M46").Value = 18200
$ws.Range("N46").Value = "$/malla 25 kilos"
$ws.Range("O46").Value = "Región de O'Higgins"
$ws.Range("P46").Value = 728
$ws.Range("Q46").Value = 25
$ws.Range("R46").Value = "Hortaliza"

# Row 47 (new)
$ws.Range("A47").Value = 1
$ws.Range("B47").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C47").Value = "Arica y Parinacota"
$ws.Range("D47").Value = 44915
$ws.Range("E47").Value = 15
$ws.Range("F47").Value = 100114001
$ws.Range("G47").Value = "Papa"
$ws.Range("H47").Value = "Asterix"
$ws.Range("I47").Value = "1a nueva(o)"
$ws.Range("J47").Value = 850
$ws.Range("K47").Value = 19000
$ws.Range("L47").Value = 20000
$ws.Range("M47").Value = 19529
$ws.Range("N47").Value = "$/saco 25 kilos"
$ws.Range("O47").Value = "Región de O'Higgins"
$ws.Range("P47").Value = 781
$ws.Range("Q47").Value = 25
$ws.Range("R47").Value = "Hortaliza"

# Row 48 (new)
$ws.Range("A48").Value = 1
$ws.Range("B48").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C48").Value = "Arica y Parinacota"
$ws.Range("D48").Value = 44915
$ws.Range("E48").Value = 15
$ws.Range("F48").Value = 100114001
$ws.Range("G48").Value = "Papa"
$ws.Range("H48").Value = "Rosara"
$ws.Range("I48").Value = "1a (cosecha)"
$ws.Range("J48").Value = 1000
$ws.Range("K48").Value = 18000
$ws.Range("L48").Value = 18000
$ws.Range("M48").Value = 18000
$ws.Range("N48").Value = "$/saco 25 kilos"
$ws.Range("O48").Value = "Región de O'Higgins"
$ws.Range("P48").Value = 720
$ws.Range("Q48").Value = 25
$ws.Range("R48").Value = "Hortaliza"
